# Hoang them chuc nang ql Dia Diem (100%)
# Set the "Ket Qua Cuoi Cung" (final result) column (E) values for the
# "3.6 Quan ly thong tin Dia Diem" (row 8), "3.10 Quan ly thong tin Mat Hang"
# (row 10) and "3.8 Quan ly phan cong Xe - Tai Xe" (row 13) rows on the
# PhanCongNganHan sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

$ws.Range("E10").Value = "100% (11/06/2010)"
$ws.Range("E8").Value = "100% (17/06/2010)"
$ws.Range("E13").Value = "100% (16/06/2010)"

# Update the view state to match the author's saved cursor position
# (scrolled so row 4 is at the top, with E13 as the active selection).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
